$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:J10").Value = 0

$ws.Range("S15").Select()
